# testarticles.xlsx edit
# 1) Drop the helper "Z41" boolean column (H): delete its header and the
#    scattered TRUE() marker cells entirely (not just clear the values).
# 2) The article-name formula cells in column C (=B#) end up showing a
#    trimmed (no trailing-space) cached value while the underlying B-column
#    text still carries its original trailing space - a pre-existing
#    data mismatch in the source data that we reproduce as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) remove column H entirely (header "Z41" + TRUE() marker cells) ---
$ws.Range("H1:H14").Clear()

# --- 2) reproduce the B/C cached-value mismatch on the name rows ---
# Excel keeps C# (=B#) perfectly in sync with B# whenever calculation is
# automatic, so to end up with a *different* cached C# value than what B#
# currently holds, we briefly go into manual calculation: set B# to the
# trimmed text (while still automatic, so C# recalculates and caches the
# trimmed text), then flip to manual and restore B#'s original
# trailing-space text. Because recalculation is suspended, C#'s cached
# value is left behind as the stale, trimmed text.

$names = @{
    2  = "Akupunkturnålar"
    3  = "Urindroppsamlare"
    8  = "Portnål"
    9  = "Syrgasmask"
    10 = "Trachealkanyl"
    12 = "Lancett"
    13 = "Provtagningskanyl"
    14 = "Blodgasspruta"
}

# Switching the Calculation property back to Automatic forces an immediate
# full recalculation, which would wipe out any stale cache already staged
# for a previous row. So: first push every row to its trimmed text while
# still automatic (each one recalculates cleanly in turn), THEN flip to
# manual a single time and restore every row's original trailing-space
# text in one batch, with no automatic/manual toggling in between.

$excel.Calculation = -4105 # xlCalculationAutomatic
foreach ($row in $names.Keys) {
    $ws.Range("B$row").Value = $names[$row]
}

$excel.Calculation = -4135 # xlCalculationManual
foreach ($row in $names.Keys) {
    $ws.Range("B$row").Value = $names[$row] + " "
}

# Leave the workbook in manual calculation so the stale C-column cache
# values created above survive the save (switching back to automatic
# would immediately resync them with column B).
